$wb = $excel.ActiveWorkbook

# --- Update metadata timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "29 Oct 2025, 06:35 PM"

# --- Update Top Gainers rows 61-76 ---
$gainers = $wb.Worksheets.Item("Top Gainers")

$data = @(
    @("GREENLAM",   3.8946, 3.5858, 10.721),
    @("APARINDS",   3.8924, 8.3414, 15.5876),
    @("HITECHGEAR", 3.8587, 1.1486, 9.9254),
    @("ORIENTTECH", 3.827,  0.5247000000000001, 32.6784),
    @("ICRA",       3.7985, 4.4793, 2.8828),
    @("SALASAR",    3.7935, 4.7872, 11.0485),
    @("NPST",       3.7841, -2.0689, -3.5677),
    @("DCW",        3.7544, 2.3219, -3.9753),
    @("RHETAN",     3.754,  4.178,  6.549),
    @("HINDPETRO",  3.6935, 6.9335, 5.7397),
    @("BHARTIHEXA", 3.6718, 7.0877, 15.3332),
    @("HLEGLAS",    3.659,  8.115500000000001, 27.1239),
    @("RHIM",       3.6544, 3.2276, 5.1826),
    @("SHK",        3.6347, 2.388,  -1.932),
    @("BCLIND",     3.6271, 2.2945, 0.1728),
    @("MUKANDLTD",  3.6133, 11.9685, 9.550800000000001)
)

$startRow = 61
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]
    $gainers.Cells.Item($row, 2).Value = $entry[0]
    $gainers.Cells.Item($row, 3).Value = $entry[1]
    $gainers.Cells.Item($row, 4).Value = $entry[2]
    $gainers.Cells.Item($row, 5).Value = $entry[3]
}
